# "Fruta / hortaliza, semanal" — weekly update.
# A new price record (Alcachofa, Argentina(o), Vega Monumental Concepción)
# is inserted as row 26, pushing the existing rows 26-51 down to 27-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting everything below it down by one.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 44757
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112013
$ws.Range("G26").Value = "Alcachofa"
$ws.Range("H26").Value = "Argentina(o)"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14400
$ws.Range("N26").Value = "$/caja 50 unidades"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 288
$ws.Range("Q26").Value = 50
$ws.Range("R26").Value = "Hortaliza"
